$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aerobox")

$ws.Range("B1").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B1").Value = 0.75
$ws.Range("B2").Value = 0.75
$ws.Range("B3").Value = 0.41666666666666669

$ws.Range("B1").Select()
